$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff": a brand-new source file
# (322329c0-2b18-49ec-9a7f-5d700d249324.md) has just been queued for
# handoff. It becomes the newest row (row 2) on every sheet, pushing the
# previous newest entry (69a0d91d-f13d-4e8e-b356-e3c5a66b8d32.md, with its
# already-completed handoff) down to row 3.
# ---------------------------------------------------------------------------

$oldBase = "69a0d91d-f13d-4e8e-b356-e3c5a66b8d32"
$newBase = "322329c0-2b18-49ec-9a7f-5d700d249324"

$oldMd = "$oldBase.md"
$newMd = "$newBase.md"

$oldHash = "94b41d1953544b5870e6c6f346a7b1dada5605ec"
$newHash = "e53757b1ad7401cffb0cfb84b7d289f536309847"

$mdUrlBase = "https://github.com/OpenLocalizationTest/oltest/blob/b44c3c8fdbf1357afbb716656531807fe46da38d/e2e"

# =================================================================
# Sheet "Overview"
# =================================================================
$ws = $wb.Worksheets.Item("Overview")

$ws.Rows(2).Insert()

# the pre-existing hyperlink geometrically stayed on row 2 - repoint it to
# row 3, where its cell data now actually lives
$ws.Range("A2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A3"), "$mdUrlBase/$oldMd", "", "", $oldMd) | Out-Null

$ws.Range("A2").Value = $newMd
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-30-21 02:30:09"
$ws.Hyperlinks.Add($ws.Range("A2"), "$mdUrlBase/$newMd", "", "", $newMd) | Out-Null

# =================================================================
# Sheet "zh-cn"
# =================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows(2).Insert()

$oldXlf = "$oldBase.$oldHash.zh-cn.xlf"
$newXlf = "$newBase.$newHash.zh-cn.xlf"
$xlfUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("D2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A3"), "$mdUrlBase/$oldMd", "", "", $oldMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "$mdUrlBase/$oldMd", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "$xlfUrlBase/2cb1cec5b118cd4371e6a97dde0215720372e6de/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldXlf", "", "", $oldXlf) | Out-Null

$ws.Range("A2").Value = $newMd
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = $newXlf
$ws.Range("E2").Value = "2016-03-21 02:30:05"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "$mdUrlBase/$newMd", "", "", $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "$mdUrlBase/$newMd", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "$xlfUrlBase/e53757b1ad7401cffb0cfb84b7d289f536309847/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlf", "", "", $newXlf) | Out-Null

# =================================================================
# Sheet "de-de"
# =================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows(2).Insert()

$oldXlfDe = "$oldBase.$oldHash.de-de.xlf"
$newXlfDe = "$newBase.$newHash.de-de.xlf"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("D2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A3"), "$mdUrlBase/$oldMd", "", "", $oldMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "$mdUrlBase/$oldMd", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "$xlfUrlBase/b89f37bf4c1a4d09de90e4e1abb1eadd41e4eb31/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldXlfDe", "", "", $oldXlfDe) | Out-Null

$ws.Range("A2").Value = $newMd
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = $newXlfDe
$ws.Range("E2").Value = "2016-03-21 02:30:09"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "$mdUrlBase/$newMd", "", "", $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "$mdUrlBase/$newMd", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "$xlfUrlBase/e53757b1ad7401cffb0cfb84b7d289f536309847/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlfDe", "", "", $newXlfDe) | Out-Null

Write-Host "Report regenerated for handoff of $newMd"
